$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase the "Nom" column (B) and fill in the generated "IDs reçus" column (D)
# for each adherent row (row 1 is the header).
$rows = @(
    @{ Row = 2;  Nom = "CARY";      Id = "231231AC0" },
    @{ Row = 3;  Nom = "TAMALOUT";  Id = "231231AT0" },
    @{ Row = 4;  Nom = "PROVOST";   Id = "231231VP0" },
    @{ Row = 5;  Nom = "PREVOST";   Id = "231231AP0" },
    @{ Row = 6;  Nom = "JULIEN";    Id = "231231MJ0" },
    @{ Row = 7;  Nom = "MAURICE";   Id = "231231MM0" },
    @{ Row = 8;  Nom = "VINCENT";   Id = "231231VV0" },
    @{ Row = 9;  Nom = "BOULANGER"; Id = "231231LB0" },
    @{ Row = 10; Nom = "PASTEUR";   Id = "231231TP0" },
    @{ Row = 11; Nom = "PASTIER";   Id = "231231JP0" },
    @{ Row = 12; Nom = "RUCHER";    Id = "231231FR0" },
    @{ Row = 13; Nom = "TISON";     Id = "231231GT0" },
    @{ Row = 14; Nom = "TRESSEAU";  Id = "231231YT0" },
    @{ Row = 15; Nom = "MORGAN";    Id = "231231GM0" },
    @{ Row = 16; Nom = "LECLERC";   Id = "231231LL0" },
    @{ Row = 17; Nom = "SOLOGNOT";  Id = "231231LS0" },
    @{ Row = 18; Nom = "ZANZIBAR";  Id = "231231MZ0" },
    @{ Row = 19; Nom = "NORMAL";    Id = "231231PN0" },
    @{ Row = 20; Nom = "FROID";     Id = "231231PF0" },
    @{ Row = 21; Nom = "CHAUD";     Id = "231231HC0" },
    @{ Row = 22; Nom = "TIÈDE";     Id = "231231AT0" },
    @{ Row = 23; Nom = "TEMPERET";  Id = "231231VT0" },
    @{ Row = 24; Nom = "PARASOL";   Id = "231231NP0" },
    @{ Row = 25; Nom = "BAIGNOIRE"; Id = "231231JB0" }
)

foreach ($entry in $rows) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Nom
    $ws.Cells.Item($entry.Row, 4).Value = $entry.Id
}

# Normalize "Jean-Pierre" -> "Jean-pierre" for the two adherents named so (rows 11 and 25)
$ws.Cells.Item(11, 3).Value = "Jean-pierre"
$ws.Cells.Item(25, 3).Value = "Jean-pierre"
